$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    @('D2', '288.98'),
    @('D3', '31.07'),
    @('E3', '2.18%'),
    @('D4', '4.955'),
    @('E4', '-0.02%'),
    @('D5', '0.07350'),
    @('E5', '1.34%'),
    @('D6', '2.343'),
    @('E6', '30.59%'),
    @('D7', '7.726'),
    @('E7', '1.76%'),
    @('B8', 'MXToken'),
    @('C8', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @('D8', '0.9115'),
    @('E8', '0.95%'),
    @('B9', 'LiechtensteinCryptoassetsExchange'),
    @('C9', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D9', '0.09295'),
    @('E9', '19.22%'),
    @('B10', 'WazirX'),
    @('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D10', '0.1706'),
    @('E10', '2.42%'),
    @('B11', 'MandalaExchangeToken'),
    @('C11', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D11', '0.08176'),
    @('E11', '2.48%'),
    @('B12', 'BitrueCoin'),
    @('C12', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D12', '0.03112'),
    @('E12', '2.01%'),
    @('B13', 'BitMartToken'),
    @('C13', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D13', '0.09980'),
    @('E13', '-0.37%'),
    @('B14', 'BitForexToken'),
    @('C14', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D14', '0.001497'),
    @('E14', '-0.29%'),
    @('B15', 'TigerCash'),
    @('C15', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @('D15', '0.005705'),
    @('E15', '-2.17%'),
    @('B16', 'LEO'),
    @('C16', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('D16', '3.470'),
    @('E16', '-0.01%'),
    @('B17', 'GateToken'),
    @('C17', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'),
    @('D17', '3.726'),
    @('E17', '0.38%'),
    @('D18', '2.091'),
    @('E18', '1.11%'),
    @('E19', '0.66%'),
    @('D20', '0.1287'),
    @('E20', '-0.56%'),
    @('D21', '4.163'),
    @('E21', '4.76%'),
    @('E22', '0.10%'),
    @('D23', '0.04521'),
    @('E23', '0.54%'),
    @('E24', '-0.14%'),
    @('D25', '0.004176'),
    @('E25', '-9.79%'),
    @('D26', '0.0001301'),
    @('E26', '0.06%'),
    @('D27', '0.0003394'),
    @('D39', '0.01582'),
    @('E39', '1.21%'),
    @('D40', '0.04473'),
    @('E40', '3.25%'),
    @('D41', '0.007365'),
    @('E41', '0.64%'),
    @('D42', '0.009882'),
    @('E42', '-1.79%'),
    @('E43', '1.98%'),
    @('D44', '0.002241'),
    @('E44', '11.44%'),
    @('D45', '0.008787'),
    @('E45', '-6.84%'),
    @('E46', '3.32%'),
    @('E47', '0.06%'),
    @('D48', '2.565'),
    @('E48', '13.77%'),
    @('D49', '0.002000'),
    @('E49', '-33.32%'),
    @('D50', '0.00002101'),
    @('E50', '0.06%'),
    @('D51', '0.0002001'),
    @('E51', '0.06%'),
)

foreach ($pair in $cells) {
    $ref = $pair[0]
    $val = $pair[1]
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}
